$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.515.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.73%  "

$ws.Range("D3").Value = "'2.648.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.84%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'510.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.20%  "

$ws.Range("D6").Value = "'157.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.91%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").Value = "'2.675.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.23%  "

$ws.Range("D10").Value = "'6.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.95%  "

$ws.Range("E11").Value = "  +5.19%  "

$ws.Range("D12").Value = "'0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.82%  "

$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("D14").Value = "'3.130.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.69%  "

$ws.Range("D15").Value = "'60.736.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.71%  "

$ws.Range("D16").Value = "'21.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.40%  "

$ws.Range("E17").Value = "  +4.96%  "

$ws.Range("E18").Value = "  +10.10%  "

$ws.Range("D19").Value = "'4.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").Value = "'351.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.34%  "

$ws.Range("D21").Value = "'10.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.50%  "

$ws.Range("D22").Value = "'6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.77%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "'60.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.78%  "

$ws.Range("D25").Value = "'0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.41%  "

$ws.Range("D26").Value = "'2.778.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.18%  "

$ws.Range("E27").Value = "  +4.86%  "

$ws.Range("D28").Value = "'0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("D29").Value = "'0.0₃0866"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.80%  "

$ws.Range("D30").Value = "'7.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.88%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'157.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.80%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'19.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.21%  "

$ws.Range("E34").Value = "  +3.91%  "

$ws.Range("D35").Value = "'5.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.01%  "

$ws.Range("D36").Value = "'4.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.91%  "

$ws.Range("E37").Value = "  +6.27%  "

$ws.Range("D38").Value = "'1.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.51%  "

$ws.Range("D39").Value = "'0.863"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.22%  "

$ws.Range("D40").Value = "'307.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.81%  "

$ws.Range("D41").Value = "'3.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.09%  "

$ws.Range("D42").Value = "'0.835"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +28.88%  "

$ws.Range("D43").Value = "'35.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.99%  "

$ws.Range("D44").Value = "'0.647"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.33%  "

$ws.Range("D45").Value = "'0.0575"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.38%  "

$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'20.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.60%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'0.994"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("D49").Value = "'4.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.17%  "

$ws.Range("E50").Value = "  +4.11%  "

$ws.Range("D51").Value = "'2.036.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.39%  "
